$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Course code text updated (M358K section header)
$ws.Range("P1").Value = "SDS348"

# New readings / micro section data added
$ws.Range("AB6").Value = 91
$ws.Range("AD6").Value = 83.7

$ws.Range("Q9").Value = 20
$ws.Range("I10").Value = 12
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 50
$ws.Range("Q25").Value = 10
$ws.Range("Q26").Value = 10
$ws.Range("Q36").Value = 5

# R36 becomes a plain literal value (no longer a formula)
$ws.Range("R36").Value = 5

# Final score formula now includes one more dropped-quiz / project column
$ws.Range("N29").Formula = '=AVERAGE(L3:L8)*M3+AVERAGE(L11:L16)*M11+SUMPRODUCT(K21:K27,M21:M27)'

# New LARGE() formulas for the 8th-ranked entries in each quiz/project block
$ws.Range("L10").Formula = '=LARGE($K$3:$K$10,8)'
$ws.Range("L18").Formula = '=LARGE($K$11:$K$20,8)'

# Update the active selection to match where the author was last working
$ws.Range("AD7").Select()
